$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the helper lookup table (N16:Q18 header/result area, N20:O32 TM/MASA MINIMA table)
$ws.Range("N16:Q32").Value = ""

# Clear the pre-filled UND column placeholder values that are being "reinforced" (left blank for user input)
$ws.Range("I31").Value = ""
$ws.Range("I37:I39").Value = ""

# Clear pre-filled equipment codes so the dropdown starts blank
$ws.Range("J42").Value = ""
$ws.Range("J43").Value = ""
$ws.Range("J45").Value = ""
